$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds "Price" text that can look like a plain number (e.g. "0.999").
# The source file stores these as literal text (inlineStr), so writing them
# straight to .Value would let Excel's type inference coerce them into real
# numbers. Force the cell to Text format before writing, then clear the
# format back afterwards so no stray number-format style is left behind on
# cells that originally had the default style.
function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue "D2" "74.346.34"
$ws.Range("E2").Value = "  +6.77%  "

Set-TextValue "D3" "2.646.45"
$ws.Range("E3").Value = "  +8.25%  "

Set-TextValue "D4" "0.999"
$ws.Range("E4").Value = "  -0.11%  "

Set-TextValue "D5" "186.02"
$ws.Range("E5").Value = "  +12.15%  "

Set-TextValue "D6" "583.70"
$ws.Range("E6").Value = "  +3.35%  "

$ws.Range("E7").Value = "  -0.18%  "

Set-TextValue "D8" "0.531"
$ws.Range("E8").Value = "  +3.88%  "

$ws.Range("E9").Value = "  +9.24%  "

Set-TextValue "D10" "2.646.84"
$ws.Range("E10").Value = "  +8.36%  "

$ws.Range("E11").Value = "  +0.77%  "

$ws.Range("E12").Value = "  +5.77%  "

$ws.Range("E13").Value = "  +0.27%  "

Set-TextValue "D14" "3.151.85"
$ws.Range("E14").Value = "  +8.91%  "

Set-TextValue "D15" "74.292.31"
$ws.Range("E15").Value = "  +6.92%  "

$ws.Range("E16").Value = "  +2.14%  "

Set-TextValue "D17" "26.19"
$ws.Range("E17").Value = "  +9.36%  "

Set-TextValue "D18" "2.655.24"
$ws.Range("E18").Value = "  +8.78%  "

Set-TextValue "D19" "9.28"
$ws.Range("E19").Value = "  +30.87%  "

Set-TextValue "D20" "11.86"
$ws.Range("E20").Value = "  +10.00%  "

Set-TextValue "D21" "368.82"
$ws.Range("E21").Value = "  +8.02%  "

Set-TextValue "D22" "2.26"
$ws.Range("E22").Value = "  +12.50%  "

Set-TextValue "D23" "4.05"
$ws.Range("E23").Value = "  +4.60%  "

Set-TextValue "D24" "6.20"
$ws.Range("E24").Value = "  +2.13%  "

Set-TextValue "D25" "0.999"
$ws.Range("E25").Value = "  -0.16%  "

Set-TextValue "D26" "69.86"
$ws.Range("E26").Value = "  +5.67%  "

Set-TextValue "D27" "4.08"
$ws.Range("E27").Value = "  +5.43%  "

Set-TextValue "D28" "9.24"
$ws.Range("E28").Value = "  +8.51%  "

Set-TextValue "D29" "2.786.92"
$ws.Range("E29").Value = "  +8.42%  "

Set-TextValue "D30" "0.998"
$ws.Range("E30").Value = "  +0.22%  "

Set-TextValue "D31" "0.0₃0928"
$ws.Range("E31").Value = "  +8.76%  "

Set-TextValue "D32" "515.55"
$ws.Range("E32").Value = "  +15.18%  "

Set-TextValue "D33" "1.38"
$ws.Range("E33").Value = "  +11.37%  "

Set-TextValue "D34" "7.60"
$ws.Range("E34").Value = "  +4.11%  "

$ws.Range("E35").Value = "  +7.02%  "

$ws.Range("E36").Value = "  -0.16%  "

Set-TextValue "D37" "162.69"
$ws.Range("E37").Value = "  +0.47%  "

Set-TextValue "D38" "0.118"
$ws.Range("E38").Value = "  +7.89%  "

Set-TextValue "D39" "19.10"
$ws.Range("E39").Value = "  +5.45%  "

Set-TextValue "D40" "19.30"
$ws.Range("E40").Value = "  +1.27%  "

$ws.Range("E41").Value = "  +0.08%  "

Set-TextValue "D42" "4.88"
$ws.Range("E42").Value = "  +10.07%  "

# Rows 43 and 44 swap coins: PolygonEcosystemToken <-> Aave.
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D43" "163.96"
$ws.Range("E43").Value = "  +24.52%  "

$ws.Range("B44").Value = "PolygonEcosystemToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-TextValue "D44" "0.325"
$ws.Range("E44").Value = "  +6.47%  "

Set-TextValue "D45" "1.65"
$ws.Range("E45").Value = "  +7.34%  "

Set-TextValue "D46" "1.18"
$ws.Range("E46").Value = "  +8.74%  "

Set-TextValue "D47" "38.99"
$ws.Range("E47").Value = "  +3.52%  "

$ws.Range("E48").Value = "  +8.92%  "

Set-TextValue "D49" "0.0851"
$ws.Range("E49").Value = "  +17.40%  "

$ws.Range("E50").Value = "  +6.18%  "

Set-TextValue "D51" "0.524"
$ws.Range("E51").Value = "  +7.23%  "
